$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-18.04479431438038, 2.398344508907384, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038),
    @(-18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038),
    @(-18.04479431438038, 2.046714678009531, 2.79429642035929, -18.04479431438038, 2.477113746341906, -18.04479431438038, 1.618584845529816, -18.04479431438038, 2.48744459571698, -18.04479431438038),
    @(-18.04479431438038, 0.8433692147553958, -18.04479431438038, -18.04479431438038, -18.04479431438038, 2.082173512163548, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038),
    @(-18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038),
    @(2.978765813982905, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038),
    @(-18.04479431438038, -18.04479431438038, -18.04479431438038, 2.837467361146716, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038),
    @(3.598939033276082, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038),
    @(-18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, 2.18332093778022),
    @(-18.04479431438038, -18.04479431438038, -18.04479431438038, 1.898414442293085, -18.04479431438038, 2.166615165982682, -18.04479431438038, -18.04479431438038, -18.04479431438038, 1.299113556962197),
    @(-18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038),
    @(-18.04479431438038, -18.04479431438038, -18.04479431438038, 1.54615267730856, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, 2.370201564496607, 1.438110124391761),
    @(-18.04479431438038, -18.04479431438038, 1.590514632012957, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, 1.956395314951585),
    @(-18.04479431438038, -18.04479431438038, -0.3692989449581339, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038),
    @(-18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, 2.557760484140287, -18.04479431438038),
    @(-18.04479431438038, 0.4299771153157192, -0.2725506293387193, -18.04479431438038, -18.04479431438038, -18.04479431438038, -0.09088776283053976, -18.04479431438038, 0.6924167125473182, -18.04479431438038),
    @(-18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, -18.04479431438038, 0.196001675234807, 4.321923026444524, 0.7807763134147587, -18.04479431438038),
    @(-18.04479431438038, -18.04479431438038, 1.809339011302504, -18.04479431438038, -18.04479431438038, -18.04479431438038, 1.725090329835201, -18.04479431438038, -18.04479431438038, -18.04479431438038),
    @(-18.04479431438038, 1.811377380402388, 2.305861121092053, -18.04479431438038, 3.851210224831434, -18.04479431438038, 2.544268840137302, -18.04479431438038, -18.04479431438038, 2.679563820577958),
    @(-18.04479431438038, 1.980065337302194, -18.04479431438038, 2.633173881466267, -18.04479431438038, 3.495158194660358, 2.512382346895391, -18.04479431438038, -18.04479431438038, -18.04479431438038)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowData = $data[$r]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowData[$c]
    }
}
